$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. In the source sheet
# this shows up as a single row inserted right above the old row 107,
# pushing every subsequent record (old rows 107-203) down by one (to rows
# 108-204), and the dimension grows from A1:R203 to A1:R204.

# Insert a new row above row 107 (shifts rows 107:203 down to 108:204).
$ws.Rows(107).Insert()

# The new row reuses the same market / category / quality / unit / origin
# metadata as the record that is now sitting in row 108 (those fields do not
# change from one weekly record to the next for this series) - only the
# date, volume, prices and $/Kg differ.
$ws.Range("A107").Value = $ws.Range("A108").Value2
$ws.Range("B107").Value = $ws.Range("B108").Value2
$ws.Range("C107").Value = $ws.Range("C108").Value2
$ws.Range("E107").Value = $ws.Range("E108").Value2
$ws.Range("F107").Value = $ws.Range("F108").Value2
$ws.Range("G107").Value = $ws.Range("G108").Value2
$ws.Range("H107").Value = $ws.Range("H108").Value2
$ws.Range("I107").Value = $ws.Range("I108").Value2
$ws.Range("N107").Value = $ws.Range("N108").Value2
$ws.Range("O107").Value = $ws.Range("O108").Value2
$ws.Range("Q107").Value = $ws.Range("Q108").Value2
$ws.Range("R107").Value = $ws.Range("R108").Value2

# New record's own values: Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg
$ws.Range("D107").Value = 44512
$ws.Range("J107").Value = 50
$ws.Range("K107").Value = 10000
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 10000
$ws.Range("P107").Value = 167
